$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The title paragraph goes from:
#   "Reflective Journal 4 (December) "
# to:
#   "Reflective Journal 5 (January) "
# with the "_GoBack" bookmark moving from the end of the "Month: Jan 2017"
# paragraph to right after "Reflective Journal 5" in the title.
# ---------------------------------------------------------------------------

# 1. Remove the old "_GoBack" bookmark (currently sits right after "2017").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Update the title text itself (scoped to paragraph 1 so nothing else in
#    the document is touched).
$p1 = $d.Paragraphs(1)
$p1.Range.Find.Execute("Reflective Journal 4 (December)", $true, $false, $false, $false, $false, $true, 1, $false, "Reflective Journal 5 (January)", 2)

# 3. Work out the character offsets (within paragraph 1) of each piece of
#    text we need as its own run: "Reflective Journal 5" | " (" | "January" | ")" | " "
$p1 = $d.Paragraphs(1)
$base = $p1.Range.Start
$text = $p1.Range.Text

$idxAfterJournal5   = $base + $text.IndexOf("Reflective Journal 5") + "Reflective Journal 5".Length
$idxAfterOpenParen  = $base + $text.IndexOf(" (") + " (".Length
$idxAfterJanuary    = $base + $text.IndexOf("January") + "January".Length
$idxAfterCloseParen = $base + $text.IndexOf(")") + ")".Length

# 4. Drop temporary bookmarks at each split point: an intervening bookmark
#    stops adjoining same-formatted runs from being coalesced back together.
$d.Bookmarks.Add("ZZZ_SPLIT1", $d.Range($idxAfterOpenParen, $idxAfterOpenParen))
$d.Bookmarks.Add("ZZZ_SPLIT2", $d.Range($idxAfterJanuary, $idxAfterJanuary))
$d.Bookmarks.Add("ZZZ_SPLIT3", $d.Range($idxAfterCloseParen, $idxAfterCloseParen))

# 5. Re-create "_GoBack" right after "Reflective Journal 5" (its new home).
$d.Bookmarks.Add("_GoBack", $d.Range($idxAfterJournal5, $idxAfterJournal5))

# 6. Touch each fragment (still scoped to paragraph 1, one replace only) so
#    the saved XML reports xml:space="preserve" exactly where it is needed.
$p1 = $d.Paragraphs(1)
$p1.Range.Find.Execute("Reflective Journal 5", $true, $false, $false, $false, $false, $true, 1, $false, "Reflective Journal 5", 1)
$p1 = $d.Paragraphs(1)
$p1.Range.Find.Execute("January", $true, $false, $false, $false, $false, $true, 1, $false, "January", 1)
$p1 = $d.Paragraphs(1)
$p1.Range.Find.Execute(")", $true, $false, $false, $false, $false, $true, 1, $false, ")", 1)

# 7. Drop the temporary split bookmarks - the run boundaries they created
#    remain in place even once the bookmarks themselves are gone.
$d.Bookmarks("ZZZ_SPLIT1").Delete()
$d.Bookmarks("ZZZ_SPLIT2").Delete()
$d.Bookmarks("ZZZ_SPLIT3").Delete()

$p1 = $d.Paragraphs(1)
Write-Host "Title now reads: [$($p1.Range.Text)]"
